$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ad"
$ws.Range("F9").Value = "ldm"
$ws.Range("F9").Select() | Out-Null
